# TC05_Bento_MultiFilter_Diagnosis-Recurrence-TumorSize-Chemo-ERStatus
# "updated bento tc as per bento perf data availability"
#
# The 4 Neo4j/Cypher queries stored on the "startup" sheet (columns B and C,
# rows 2-4) all filter on tp.chemotherapy_regimen. The chemotherapy regimen
# filter value is being swapped from
#   "Other treatment given as part of a CTSU protocol"
# to
#   "Dose dense AC (2 week cycles)"
# (note the extra space left behind before "and d.er_status" in the new text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSnippet = 'tp.chemotherapy_regimen In ["Other treatment given as part of a CTSU protocol"] and d.er_status'
$newSnippet = 'tp.chemotherapy_regimen In ["Dose dense AC (2 week cycles)"]  and d.er_status'

foreach ($addr in @("B2", "C2", "B3", "C3", "B4", "C4")) {
    $cell = $ws.Range($addr)
    $current = $cell.Value2
    if ($current -ne $null -and $current.Contains($oldSnippet)) {
        $cell.Value2 = $current.Replace($oldSnippet, $newSnippet)
    }
}

# Reflect the author's final cursor position/selection on the sheet.
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
